$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 4679.8
$ws.Range("I28").Value = 900
$ws.Range("K28").Value = 900
$ws.Range("M28").Value = -415

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2624.75
$ws.Range("J40").Value = 2999.6667
$ws.Range("L40").Value = 2999.6667
$ws.Range("N40").Value = -3349.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 3201.4443
$ws.Range("I107").Value = 833.6667
$ws.Range("K107").Value = 833.6667
$ws.Range("M107").Value = 1086.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6000
$ws.Range("J116").Value = 6000
$ws.Range("L116").Value = 6000
$ws.Range("N116").Value = -12884

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2677.9
$ws.Range("I137").Value = 2197.75
$ws.Range("J137").Value = 2998
$ws.Range("K137").Value = 6593.25
$ws.Range("L137").Value = 8994
$ws.Range("M137").Value = -4043.25
$ws.Range("N137").Value = -14094

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4172.143
$ws.Range("I141").Value = 3724.077
$ws.Range("K141").Value = 11172.231
$ws.Range("M141").Value = -5992.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 1389.2
$ws.Range("I35").Value = 700
$ws.Range("K35").Value = 700
$ws.Range("M35").Value = -294

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3303.5
$ws.Range("I61").Value = 3303.5
$ws.Range("K61").Value = 3303.5
$ws.Range("M61").Value = -3091.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3566.3333
$ws.Range("I63").Value = 2699
$ws.Range("K63").Value = 2699
$ws.Range("M63").Value = -2013

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3566.3333
$ws.Range("I66").Value = 2699
$ws.Range("K66").Value = 13495
$ws.Range("M66").Value = -10063

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2270.7778
$ws.Range("I88").Value = 863.1111
$ws.Range("K88").Value = 863.1111
$ws.Range("M88").Value = -457.1111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2270.7778
$ws.Range("I91").Value = 863.1111
$ws.Range("K91").Value = 863.1111
$ws.Range("M91").Value = 540.8889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 44250
$ws.Range("J125").Value = 44250
$ws.Range("L125").Value = 44250
$ws.Range("N125").Value = -54090

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1962.7858
$ws.Range("I132").Value = 770
$ws.Range("K132").Value = 2310
$ws.Range("M132").Value = 220

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3303.5
$ws.Range("I136").Value = 3303.5
$ws.Range("K136").Value = 9910.5
$ws.Range("M136").Value = -7360.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 850
$ws.Range("I37").Value = 300
$ws.Range("J37").Value = 1400
$ws.Range("K37").Value = 300
$ws.Range("L37").Value = 1400
$ws.Range("M37").Value = -163
$ws.Range("N37").Value = -1674

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5409.0454
$ws.Range("I134").Value = 5409.0454
$ws.Range("K134").Value = 16227.1362
$ws.Range("M134").Value = -13692.1362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 68.59999999999999
$ws.Range("I7").Value = 68.59999999999999
$ws.Range("K7").Value = 68.59999999999999
$ws.Range("M7").Value = 44.40000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7941.25
$ws.Range("I99").Value = 7941.25
$ws.Range("K99").Value = 7941.25
$ws.Range("M99").Value = -6443.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 7941.25
$ws.Range("I126").Value = 7941.25
$ws.Range("K126").Value = 23823.75
$ws.Range("M126").Value = -21353.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2391.25
$ws.Range("I132").Value = 1440.75
$ws.Range("K132").Value = 4322.25
$ws.Range("M132").Value = -1792.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4180.923
$ws.Range("I134").Value = 4131.4
$ws.Range("K134").Value = 12394.2
$ws.Range("M134").Value = -9859.199999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2019.375
$ws.Range("I75").Value = 2482
$ws.Range("J75").Value = 1741.8
$ws.Range("K75").Value = 7446
$ws.Range("L75").Value = 5225.4
$ws.Range("M75").Value = -6448
$ws.Range("N75").Value = -7221.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 2019.375
$ws.Range("I78").Value = 2482
$ws.Range("J78").Value = 1741.8
$ws.Range("K78").Value = 22338
$ws.Range("L78").Value = 15676.2
$ws.Range("M78").Value = -17346
$ws.Range("N78").Value = -25660.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 224
$ws.Range("I103").Value = 223.66667
$ws.Range("J103").Value = 225
$ws.Range("K103").Value = 671.00001
$ws.Range("L103").Value = 675
$ws.Range("M103").Value = 207.99999
$ws.Range("N103").Value = -2433

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 502
$ws.Range("I108").Value = 502
$ws.Range("K108").Value = 1506
$ws.Range("M108").Value = 1374

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 887.2
$ws.Range("I114").Value = 609.5
$ws.Range("J114").Value = 1072.3334
$ws.Range("K114").Value = 1828.5
$ws.Range("L114").Value = 3217.0002
$ws.Range("M114").Value = 1425.5
$ws.Range("N114").Value = -9725.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 871.75
$ws.Range("I117").Value = 768.5
$ws.Range("K117").Value = 2305.5
$ws.Range("M117").Value = 1136.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 922.82355
$ws.Range("I121").Value = 701.625
$ws.Range("J121").Value = 1119.4445
$ws.Range("K121").Value = 2104.875
$ws.Range("L121").Value = 3358.3335
$ws.Range("M121").Value = -794.875
$ws.Range("N121").Value = -5978.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 705
$ws.Range("I122").Value = 765.5
$ws.Range("K122").Value = 6889.5
$ws.Range("M122").Value = -4439.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 625
$ws.Range("I129").Value = 656.25
$ws.Range("K129").Value = 1968.75
$ws.Range("M129").Value = 3031.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1438.625
$ws.Range("I131").Value = 996.25
$ws.Range("J131").Value = 1881
$ws.Range("K131").Value = 2988.75
$ws.Range("L131").Value = 5643
$ws.Range("M131").Value = 2051.25
$ws.Range("N131").Value = -15723

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 30000
$ws.Range("J52").Value = 30000
$ws.Range("L52").Value = 30000
$ws.Range("N52").Value = -30518

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4969.375
$ws.Range("I113").Value = 4352.2
$ws.Range("K113").Value = 4352.2
$ws.Range("M113").Value = -2182.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1492.6364
$ws.Range("I132").Value = 1492.6364
$ws.Range("K132").Value = 4477.9092
$ws.Range("M132").Value = -1947.9092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1675

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1675

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3333.6667
$ws.Range("I40").Value = 2748.5
$ws.Range("K40").Value = 2748.5
$ws.Range("M40").Value = -2612.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6349.75
$ws.Range("I132").Value = 5282.8335
$ws.Range("K132").Value = 15848.5005
$ws.Range("M132").Value = -13318.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8000
$ws.Range("J14").Value = 8000
$ws.Range("L14").Value = 8000
$ws.Range("N14").Value = -8336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2000
$ws.Range("I81").Value = 2000
$ws.Range("K81").Value = 4000
$ws.Range("M81").Value = -2939

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2000
$ws.Range("I84").Value = 2000
$ws.Range("K84").Value = 20000
$ws.Range("M84").Value = -14696
